$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 109
$ws.Range("A109").Value() = 'Deliberative Choice Strategies in Youths: Relevance to Transdiagnostic Anxiety Symptoms'
$ws.Range("B109").Value() = 'Elise M. Cardinale, David Pagliaccio, Caroline Swetlitz, Hannah Grassie, Rany Abend, Vincent Costa, Bruno Averbeck, Melissa A. Brotman, Daniel S. Pine, Ellen Leibenluft, Katharina Kircanski'
$ws.Range("C109").Value() = 'Clinical Psychological Science'
$ws.Range("D109").Value() = 2021
$ws.Range("E109").Value() = 'https://doi.org/10.1177%2F2167702621991805'
$ws.Range("F109").Value() = 'Aberrant decision-making characterizes various pediatric psychopathologies; however, deliberative choice strategies
have not been investigated. A transdiagnostic sample of 95 youths completed a child-friendly sequential sampling
paradigm. Participants searched for the best offer by sampling a finite list of offers. Participants’ willingness to explore
was measured as the number of offers sampled, and ideal task performance was modeled using a Markov decisionprocess model. As in previous findings in adults, youths explored more offers when lists were long compared with
short, yet participants generally sampled fewer offers relative to model-estimated ideal performance. Searching deeper
into the list was associated with choosing better price options. Analyses examining the main and interactive effects of
transdiagnostic anxiety and irritability symptoms indicated a negative correlation between anxiety and task performance
(p = .01, ηp
2 = .08). Findings suggest the need for more research on exploratory decision impairments in youths with
anxiety symptoms.'
$ws.Range("G109").Value() = 'ARI used as part of measures'
$ws.Range("H109").Value() = 'Observational'

# Row 110
$ws.Range("A110").Value() = 'Attention bias to negative versus non-negative faces is related to negative affectivity in a transdiagnostic youth sample'
$ws.Range("B110").Value() = 'Anita Harrewijn, Rany Abend, Reut Naim, Simone P. Haller, Caitlin M. Stavish, Mira A. Bajaj, Chika Matsumoto, Kelly Dombek, Elise M. Cardinale, Katharina Kircanski, Melissa A. Brotman'
$ws.Range("C110").Value() = 'Journal of Psychiatric Research'
$ws.Range("D110").Value() = 2021
$ws.Range("E110").Value() = 'https://doi.org/10.1016/j.jpsychires.2021.04.036'
$ws.Range("F110").Value() = 'This study identified a shared pathophysiological mechanism of pediatric anxiety and irritability. Clinically, anxiety and irritability are common, co-occurring problems, both characterized by high-arousal negative affective states. Behaviorally, anxiety and irritability are associated with aberrant threat processing. To build on these findings, we examined eye-tracking measures of attention bias in relation to the unique and shared features of anxiety and irritability in a transdiagnostic sample of youth (n = 97, 58% female, Mage = 13.03, SDage = 2.82). We measured attention bias to negative versus non-negative emotional faces during a passive viewing task. We employed bifactor analysis to parse the unique and shared variance of anxiety and irritability symptoms from self- and parent-report questionnaires. Negative affectivity is the derived latent factor reflecting shared variance of anxiety and irritability. We found that higher negative affectivity was associated with looking longer at negative versus non-negative faces, reflecting a shared mechanism of anxiety and irritability. This finding suggests that modification of elevated attention to negative emotional faces may represent a common potential treatment target of anxiety and irritability.'
$ws.Range("G110").Value() = 'ARI used as part of measures'
$ws.Range("H110").Value() = 'Observational'

# Row 111
$ws.Range("A111").Value() = 'Reduced neural responsiveness to looming stimuli is associated with increased aggression'
$ws.Range("B111").Value() = 'R James Blair, Ru Zhang, Johannah Bashford-Largo, Sahil Bajaj, Avantika Mathur, Jay Ringle, Amanda Schwartz, Jaimie Elowsky, Matthew Dobbertin, Karina S Blair, Patrick M Tyler'
$ws.Range("C111").Value() = 'Social Cognitive and Affective Neuroscience'
$ws.Range("D111").Value() = 2021
$ws.Range("E111").Value() = 'https://doi.org/10.1093/scan/nsab058'
$ws.Range("F111").Value() = 'While neuro-cognitive work examining aggression has examined patients with conditions at increased risk for aggression or individuals self-reporting past aggression, little work has attempted to identify neuro-cognitive markers associated with observed/recorded aggression. The goal of the current study was to determine the extent to which aggression by youth in the first three months of residential care was associated with atypical responsiveness to threat stimuli. This functional MRI study involved 98 (68 male; mean age = 15.96 [sd = 1.52]) adolescents in residential care performing a looming threat task involving images of threatening and neutral human faces or animals that appeared to be either loom or recede. Level of aggression was negatively associated with responding to looming stimuli (irrespective of whether these were threatening or neutral) within regions including bilateral inferior frontal gyrus, right inferior parietal lobule, right superior/middle temporal gyrus and a region of right uncus proximal to the amygdala. These data indicate that aggression level is associated with a decrease in responsiveness to a basic threat cue-looming stimuli. Reduced threat responsiveness likely results in the individual being less able to represent the negative consequences that may result from engaging in aggression, thereby increasing the risk for aggressive episodes.'
$ws.Range("G111").Value() = 'ARI used as part of measures'
$ws.Range("H111").Value() = 'Observational'

# Row 112
$ws.Range("A112").Value() = 'Trauma Reactions in Children with ADHD During the COVID-19 Pandemic: The Mediating Effect of Irritability'
$ws.Range("B112").Value() = 'Halit Necmi Uçar, Fatih Hilmi Çetin, Serhat Türkoğlu, Gökçeçiçek Arıcı Sağliyan, Özlem Çiçek Zekey & Çiğdem Yılmaz'
$ws.Range("C112").Value() = 'Journal of Loss and Trauma'
$ws.Range("D112").Value() = 2021
$ws.Range("E112").Value() = 'https://doi.org/10.1080/15325024.2021.1926782'
$ws.Range("F112").Value() = 'This study aimed to investigate the symptoms of trauma, depression, anxiety, and irritability in children with attention deficit hyperactivity disorder (ADHD) during home confinement due to COVID-19. The sample of this single-center, cross-sectional study consisted of 42 children with ADHD and their parents. There were significant differences in the severity of depression, anxiety, irritability, and ADHD symptoms between the group that experienced pathological trauma reactions and the group that did not. According to the results of mediation analyses, the irritability score was found to be the full mediating factor in the relationship between trauma reaction scores/emotional symptom scores and the severity of ADHD symptoms. Because of the modest sample size due to the exclusion of ADHD patients with comorbidities, we recommend that future research focus on reproducing the role of irritability, in larger samples, as the full mediator on the path from cognition to behavior in children with ADHD under home confinement during the pandemic.'
$ws.Range("G112").Value() = 'ARI used as part of measures'
$ws.Range("H112").Value() = 'Observational'

# Row 113
$ws.Range("A113").Value() = 'Bullying Perpetration and Victimization in Youth: Associations with Irritability and Anxiety'
$ws.Range("B113").Value() = 'Hung-Wei Bernie Chen, Erin S. Gardner, Tessa Clarkson, Nicholas R. Eaton, Jillian Lee Wiggins, Ellen Leibenluft & Johanna M. Jarcho'
$ws.Range("C113").Value() = 'Child Psychiatry & Human Development'
$ws.Range("D113").Value() = 2021
$ws.Range("E113").Value() = 'https://doi.org/10.1007/s10578-021-01192-9'
$ws.Range("F113").Value() = 'Prior work on has demonstrated that irritability and anxiety are associated with bullying perpetration and victimization, respectively. Even though symptoms of irritability and anxiety often occur concurrently, few studies have tested their interactive effects on perpetration or victimization. The current study recruited 131 youths from a broader program of research that examines the pathophysiology and treatment of pediatric irritability and anxiety. Two moderation tests were performed to examine concurrent irritability and anxiety symptoms and their relation to perpetration and victimization of bullying. More severe anxiety was associated with greater victimization. However, more severe irritability was associated with, not just greater perpetration, but also greater victimization. An irritability-by-anxiety interaction demonstrated that youths with more severe irritability and lower levels of anxiety engaged in more perpetration. Our findings suggest a more nuanced approach to understanding how the commonly comorbid symptoms of irritability and anxiety interact in relation to peer-directed behavior in youths.
'
$ws.Range("G113").Value() = 'ARI used as part of measures'
$ws.Range("H113").Value() = 'Observational'

# Row 114
$ws.Range("A114").Value() = 'Early Predictors of Adolescent Irritability'
$ws.Range("B114").Value() = 'Ellen M. Kessel, Lea R. Dougherty, Samantha Hubacheck, Emma Chad-Friedman, Tom Olino, Gabrielle A. Carlson, Daniel N. Klein,'
$ws.Range("C114").Value() = 'Child and Adolescent Psychiatric Clinics of North America'
$ws.Range("D114").Value() = 2021
$ws.Range("E114").Value() = 'https://doi.org/10.1016/j.chc.2021.04.002'
$ws.Range("F114").Value() = '•
Limited research has examined early precursors/risk factors for adolescent irritability.
•
Using multiple methods and informants and a prospective longitudinal design, the authors examined the continuity of irritability from early childhood to adolescence and identified other early predictors of adolescent irritability.
•
Across both self-reports and mother-reports, the authors found evidence for continuity of irritability from ages 3 to 15.
•
The authors also found that early antecedents of adolescent irritability differ in many instances as a function of either informant or gender.
•
The results also suggest that adolescent irritability is characterized by several distinct developmental pathways from age 3 that have the potential to result in an irritable phenotype at age 15.
•
They also suggest that self-reported and mother-reported irritability may be capturing distinct underlying constructs, and both should be considered in assessments of adolescent irritability.'
$ws.Range("G114").Value() = 'ARI used as part of measures'
$ws.Range("H114").Value() = 'Observational'
$ws.Range("I114").Value() = 'Psychometric'
